# Apply updated cryptocurrency market data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text interpretation for
# purely-numeric-looking strings (so Excel does not convert them to
# floating point numbers and drop formatting like trailing zeros),
# then restore the cell style so no stray style/number-format change
# is introduced.
function Set-TextValue($rng, [string]$val) {
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = "29.881.55"
$ws.Range("E2").Value = "  +0.83%  "

# Row 3
$ws.Range("D3").Value = "1.625.29"
$ws.Range("E3").Value = "  +0.96%  "

# Row 4
$ws.Range("E4").Value = "  -0.32%  "

# Row 5
Set-TextValue $ws.Range("D5") "214.71"
$ws.Range("E5").Value = "  +0.92%  "

# Row 6
$ws.Range("E6").Value = "  +0.40%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.997"
$ws.Range("E7").Value = "  -0.33%  "

# Row 8
Set-TextValue $ws.Range("D8") "29.79"
$ws.Range("E8").Value = "  +9.95%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.258"
$ws.Range("E9").Value = "  +2.64%  "

# Row 10
$ws.Range("E10").Value = "  +1.52%  "

# Row 11
$ws.Range("E11").Value = "  +0.40%  "

# Row 12
$ws.Range("D12").Value = "1.859.17"

# Row 13
$ws.Range("D13").Value = "1.624.91"
$ws.Range("E13").Value = "  +0.84%  "

# Row 14
$ws.Range("E14").Value = "  +6.33%  "

# Row 15
$ws.Range("E15").Value = "  +4.83%  "

# Row 16
$ws.Range("D16").Value = "29.952.42"
$ws.Range("E16").Value = "  +1.01%  "

# Row 17
Set-TextValue $ws.Range("D17") "8.86"
$ws.Range("E17").Value = "  +16.55%  "

# Row 18
Set-TextValue $ws.Range("D18") "64.65"
$ws.Range("E18").Value = "  +1.73%  "

# Row 19
Set-TextValue $ws.Range("D19") "244.11"
$ws.Range("E19").Value = "  +1.53%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0706"
$ws.Range("E20").Value = "  +1.66%  "

# Row 21
$ws.Range("E21").Value = "  -0.30%  "

# Row 22
$ws.Range("E22").Value = "  +3.42%  "

# Row 23
Set-TextValue $ws.Range("D23") "9.63"
$ws.Range("E23").Value = "  +4.30%  "

# Row 24
Set-TextValue $ws.Range("D24") "2.12"
$ws.Range("E24").Value = "  +1.05%  "

# Row 25
Set-TextValue $ws.Range("D25") "157.40"
$ws.Range("E25").Value = "  +1.54%  "

# Row 26
Set-TextValue $ws.Range("D26") "15.68"
$ws.Range("E26").Value = "  +2.28%  "

# Row 27
$ws.Range("E27").Value = "  +2.45%  "

# Row 28
Set-TextValue $ws.Range("D28") "6.61"
$ws.Range("E28").Value = "  +2.79%  "

# Row 29
Set-TextValue $ws.Range("D29") "0.998"
$ws.Range("E29").Value = "  -0.32%  "

# Row 30
$ws.Range("E30").Value = "  +3.03%  "

# Row 31
$ws.Range("E31").Value = "  +5.51%  "

# Row 32
$ws.Range("E32").Value = "  +3.98%  "

# Row 33
Set-TextValue $ws.Range("D33") "3.23"
$ws.Range("E33").Value = "  +3.27%  "

# Row 34
$ws.Range("D34").Value = "1.425.77"
$ws.Range("E34").Value = "  +0.29%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.64"
$ws.Range("E35").Value = "  +6.66%  "

# Row 36
$ws.Range("E36").Value = "  -0.01%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.86"
$ws.Range("E37").Value = "  +1.64%  "

# Row 38
$ws.Range("E38").Value = "  -0.70%  "

# Row 39
$ws.Range("E39").Value = "  +2.76%  "

# Row 40
$ws.Range("E40").Value = "  +3.42%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D41") "0.0502"
$ws.Range("E41").Value = "  +1.86%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D42") "1.99"
$ws.Range("E42").Value = "  +0.32%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D43") "0.832"
$ws.Range("E43").Value = "  +4.01%  "

# Row 44
Set-TextValue $ws.Range("D44") "55.00"
$ws.Range("E44").Value = "  +0.73%  "

# Row 45
Set-TextValue $ws.Range("D45") "69.30"
$ws.Range("E45").Value = "  +5.02%  "

# Row 46
$ws.Range("E46").Value = "  +15.83%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.997"
$ws.Range("E47").Value = "  -0.35%  "

# Row 48
$ws.Range("E48").Value = "  +1.85%  "

# Row 49
$ws.Range("D49").Value = "1.766.71"
$ws.Range("E49").Value = "  +0.93%  "

# Row 50
$ws.Range("E50").Value = "  +2.32%  "

# Row 51
$ws.Range("D51").Value = "0.0₆0110"
$ws.Range("E51").Value = "  +3.35%  "

Write-Output "Done applying updates"